$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.960.49"
$ws.Range("E2").Value = "  +5.91%  "
$ws.Range("D3").Value = "3.649.12"
$ws.Range("E3").Value = "  +5.65%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'593.46"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").Value = "'195.66"
$ws.Range("E6").Value = "  +3.29%  "
$ws.Range("D7").Value = "'0.648"
$ws.Range("E7").Value = "  +2.81%  "
$ws.Range("D8").Value = "3.643.04"
$ws.Range("E8").Value = "  +5.67%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  +7.44%  "
$ws.Range("D11").Value = "'0.682"
$ws.Range("E11").Value = "  +5.28%  "
$ws.Range("D12").Value = "'58.38"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "'0.0000299"
$ws.Range("E13").Value = "  +7.50%  "
$ws.Range("E14").Value = "  +5.58%  "
$ws.Range("D15").Value = "4.233.08"
$ws.Range("E15").Value = "  +5.76%  "
$ws.Range("D16").Value = "'20.44"
$ws.Range("E16").Value = "  +7.77%  "
$ws.Range("D17").Value = "3.648.24"
$ws.Range("E17").Value = "  +5.70%  "
$ws.Range("D18").Value = "70.969.81"
$ws.Range("E18").Value = "  +6.05%  "
$ws.Range("D19").Value = "'12.80"
$ws.Range("E19").Value = "  +5.37%  "
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("E21").Value = "  +3.69%  "
$ws.Range("D22").Value = "'490.53"
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("D23").Value = "'19.37"
$ws.Range("E23").Value = "  +13.31%  "
$ws.Range("D24").Value = "'5.24"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("E25").Value = "  +2.74%  "
$ws.Range("D26").Value = "'91.69"
$ws.Range("E26").Value = "  +2.14%  "
$ws.Range("D27").Value = "'3.19"
$ws.Range("E27").Value = "  +6.27%  "
$ws.Range("E28").Value = "  +4.74%  "
$ws.Range("D29").Value = "'9.64"
$ws.Range("E29").Value = "  +6.30%  "
$ws.Range("D30").Value = "'7.94"
$ws.Range("E30").Value = "  +6.21%  "
$ws.Range("D31").Value = "'32.90"
$ws.Range("E31").Value = "  +4.97%  "
$ws.Range("D32").Value = "'0.124"
$ws.Range("E32").Value = "  +9.94%  "
$ws.Range("D33").Value = "'12.33"
$ws.Range("D34").Value = "'618.37"
$ws.Range("E34").Value = "  +3.16%  "
$ws.Range("D35").Value = "'66.56"
$ws.Range("E35").Value = "  +3.58%  "
$ws.Range("D36").Value = "'40.31"
$ws.Range("E36").Value = "  +7.39%  "
$ws.Range("D37").Value = "0.0₃0835"
$ws.Range("E37").Value = "  +10.35%  "
$ws.Range("E38").Value = "  +5.81%  "
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "'3.58"
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("D42").Value = "3.334.10"
$ws.Range("E42").Value = "  +3.67%  "
$ws.Range("D43").Value = "'3.25"
$ws.Range("E43").Value = "  +16.73%  "
$ws.Range("D44").Value = "'3.19"
$ws.Range("E44").Value = "  +8.25%  "
$ws.Range("E45").Value = "  +9.91%  "
$ws.Range("D46").Value = "'0.0459"
$ws.Range("E46").Value = "  +6.14%  "
$ws.Range("E47").Value = "  +10.90%  "
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  +0.13%  "
